$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "utilizando densidade populacional para classificar grandes cidades"
#
# Reclassify some parties' ideologia (column B, keyed by sigla_partido in
# column A):
#  - PHS           : "Extinto" -> "Direita"
#  - PMN, PP, PRB, PTB : "Centrao" -> merged into the "Centro" bucket
#
# After the reassignment, "Extinto" and "Centrao" are no longer referenced by
# any row, and every former "Centro Dem" row is relabeled "Centro".

$partiesToDireita = @("PHS")
$partiesToCentro  = @("PMN", "PP", "PRB", "PTB")

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $sigla = $ws.Cells.Item($r, 1).Text

    if ($partiesToDireita -contains $sigla) {
        $ws.Cells.Item($r, 2).Value = "Direita"
    }
    elseif ($partiesToCentro -contains $sigla) {
        $ws.Cells.Item($r, 2).Value = "Centro Dem"
    }
}

# Collapse every "Centro Dem" (the original rows plus the ones just merged
# in above) down to the new, shorter "Centro" label.
$null = $ws.Cells.Replace("Centro Dem", "Centro")
